# Applies the Flashscore odds-refresh edit described in the commit:
# "Atualizando o arquivo XLSX" - updates specific odds cells on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new numeric value, grouped by data row.
$updates = [ordered]@{
    # Row 35
    "J35" = 1.18   # was 1.2
    "K35" = 4.5   # was 4.33
    # Row 56
    "AC56" = 26   # was 29
    "AE56" = 15   # was 13
    "AF56" = 17   # was 15
    "AG56" = 11   # was 9.5
    "AH56" = 26   # was 23
    "AI56" = 17   # was 15
    "AJ56" = 21   # was 19
    "G56" = 2.5   # was 2.9
    "I56" = 2.55   # was 2.25
    "N56" = 1.44   # was 1.48
    "O56" = 2.7   # was 2.6
    "P56" = 1.22   # was 1.25
    "Q56" = 4   # was 3.75
    "T56" = 15   # was 17
    "U56" = 17   # was 19
    "V56" = 10   # was 12
    "W56" = 26   # was 34
    "X56" = 17   # was 21
    "Y56" = 21   # was 23
    # Row 65
    "AE65" = 7   # was 7.5
    "AH65" = 29   # was 34
    "G65" = 2.6   # was 2.45
    "I65" = 2.9   # was 3.1
    "J65" = 1.11   # was 1.1
    "K65" = 6.5   # was 7
    "U65" = 11   # was 10
    "V65" = 11   # was 10
    "W65" = 26   # was 23
    "X65" = 26   # was 23
    "Z65" = 6.5   # was 7
    # Row 66
    "AA66" = 7   # was 7.5
    "AE66" = 12   # was 13
    "AF66" = 26   # was 29
    "AG66" = 17   # was 19
    "AH66" = 51   # was 67
    "G66" = 1.7   # was 1.55
    "H66" = 3.7   # was 3.9
    "I66" = 5.25   # was 6.25
    "J66" = 1.07   # was 1.06
    "K66" = 9   # was 10
    "L66" = 1.36   # was 1.3
    "M66" = 3.2   # was 3.5
    "N66" = 2.08   # was 1.98
    "O66" = 1.73   # was 1.88
    "T66" = 6   # was 6.5
    "U66" = 7.5   # was 7
    "W66" = 13   # was 11
    "X66" = 15   # was 13
    "Z66" = 9   # was 10
    # Row 70
    "AG70" = 23   # was 26
    "AH70" = 101   # was 126
    "AJ70" = 51   # was 67
    "I70" = 9   # was 9.5
    "U70" = 6.5   # was 6
    "Z70" = 13   # was 12
    # Row 72
    "AB72" = 21   # was 19
    "AC72" = 81   # was 67
    "AE72" = 6   # was 7
    "AF72" = 11   # was 13
    "AG72" = 11   # was 12
    "AH72" = 26   # was 29
    "AI72" = 26   # was 29
    "G72" = 2.88   # was 2.5
    "H72" = 3   # was 3.1
    "I72" = 2.6   # was 2.9
    "J72" = 1.13   # was 1.1
    "K72" = 6   # was 7
    "L72" = 1.57   # was 1.5
    "M72" = 2.25   # was 2.5
    "N72" = 2.88   # was 2.6
    "O72" = 1.4   # was 1.48
    "P72" = 1.62   # was 1.53
    "Q72" = 2.2   # was 2.38
    "R72" = 2.25   # was 2.1
    "S72" = 1.57   # was 1.67
    "U72" = 12   # was 11
    "V72" = 12   # was 11
    "W72" = 29   # was 23
    "X72" = 29   # was 23
    "Z72" = 6   # was 6.5
    # Row 73
    "AE73" = 6   # was 5.5
    "H73" = 3.3   # was 3.25
    "L73" = 1.4   # was 1.44
    "M73" = 2.75   # was 2.63
    "N73" = 2.25   # was 2.35
    "O73" = 1.62   # was 1.57
    "Z73" = 7.5   # was 7
    # Row 126
    "AE126" = 11   # was 10
    "AF126" = 21   # was 17
    "AI126" = 34   # was 29
    "G126" = 1.91   # was 2.1
    "H126" = 3.6   # was 3.3
    "I126" = 4.1   # was 3.6
    "J126" = 1.05   # was 1.04
    "K126" = 11   # was 10
    "M126" = 3.4   # was 3.25
    "N126" = 2   # was 2.1
    "O126" = 1.8   # was 1.7
    "P126" = 1.4   # was 1.44
    "Q126" = 2.75   # was 2.63
    "U126" = 8.5   # was 9.5
    "V126" = 8.5   # was 9
    "W126" = 17   # was 19
    "X126" = 15   # was 17
    "Z126" = 9.5   # was 9
    # Row 128
    "AF128" = 12   # was 13
    "AG128" = 11   # was 12
    "AI128" = 26   # was 29
    "G128" = 2.6   # was 2.45
    "I128" = 2.7   # was 2.88
    "V128" = 11   # was 10
    "W128" = 26   # was 23
    "X128" = 26   # was 23
    # Row 131
    "AA131" = 8   # was 8.5
    "AD131" = 351   # was 301
    "AG131" = 23   # was 21
    "I131" = 8   # was 7.5
    "J131" = 1.05   # was 1.04
    "K131" = 11   # was 13
    "L131" = 1.25   # was 1.22
    "M131" = 3.75   # was 4
    "N131" = 1.85   # was 1.73
    "O131" = 1.95   # was 2.08
    "P131" = 1.33   # was 1.3
    "Q131" = 3.25   # was 3.4
    "R131" = 2   # was 1.83
    "S131" = 1.73   # was 1.83
    "T131" = 6.5   # was 7.5
    "U131" = 6.5   # was 7
    "Y131" = 29   # was 26
    "Z131" = 11   # was 12
    # Row 205
    "AB205" = 13   # was 12
    "AF205" = 11   # was 10
    "AH205" = 19   # was 17
    "AJ205" = 21   # was 23
    "G205" = 3.6   # was 3.7
    "H205" = 3.6   # was 3.5
    "L205" = 1.22   # was 1.2
    "M205" = 4   # was 4.33
    "U205" = 19   # was 21
    "V205" = 12   # was 13
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
